# Corregí error de bitacora
#
# The log entries had accidentally been written starting at row 15,
# leaving rows 6-14 blank. Move that block (A15:F17) back up so it
# follows directly after the existing data (which ends at row 5),
# landing at A6:F8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cut preserves the original cell styles (e.g. the date format applied
# to column A) instead of creating new style/number-format entries the
# way a Copy+PasteSpecial would.
$src = $ws.Range("A15:F17")
$dst = $ws.Range("A6:F8")
$src.Cut($dst) | Out-Null

# Remove the now-vacated rows entirely so the sheet's dimension/used
# range shrinks back down (matches dimension going from A1:F17 to A1:F8).
$ws.Rows("15:17").Delete() | Out-Null

# Match the new selection left behind by the move.
$ws.Range("A6:H8").Select() | Out-Null

$excel.CutCopyMode = $false
